$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Semana_7")

$ws.Range("R5").Value = 15
$ws.Range("T5").Value = 0
$ws.Range("L6").Value = 0
$ws.Range("R6").Value = 5
$ws.Range("T6").Value = 0
$ws.Range("R7").Value = 1
$ws.Range("T7").Value = 3
$ws.Range("U7").Value = 6
$ws.Range("R8").Value = 3
$ws.Range("T8").Value = 0
$ws.Range("R12").Value = 1
$ws.Range("L14").Value = 0
$ws.Range("R15").Value = 1
$ws.Range("L17").Value = 0
$ws.Range("L19").Value = 0
$ws.Range("R19").Value = 1
$ws.Range("R20").Value = 1
$ws.Range("L21").Value = 0
$ws.Range("R21").Value = 3
$ws.Range("T21").Value = 0
$ws.Range("R22").Value = 3
$ws.Range("T22").Value = 0
$ws.Range("R23").Value = 2
$ws.Range("T23").Value = 7
$ws.Range("L24").Value = 0
$ws.Range("R24").Value = 8
$ws.Range("L25").Value = 0
$ws.Range("R25").Value = 4
$ws.Range("L28").Value = 0
$ws.Range("R28").Value = 6
$ws.Range("T28").Value = 0
$ws.Range("L29").Value = 0
$ws.Range("R29").Value = 11
$ws.Range("T29").Value = 9
$ws.Range("U29").Value = 16
$ws.Range("L33").Value = 0
$ws.Range("R33").Value = 3
$ws.Range("T33").Value = 0
$ws.Range("L35").Value = 0
$ws.Range("R35").Value = 8
$ws.Range("L36").Value = 0
$ws.Range("L37").Value = 0
$ws.Range("R37").Value = 5
$ws.Range("T37").Value = 0
$ws.Range("L40").Value = 0
$ws.Range("R42").Value = 13
$ws.Range("T42").Value = 0
$ws.Range("L44").Value = 0
$ws.Range("R44").Value = 2
$ws.Range("L46").Value = 0
$ws.Range("R46").Value = 1
$ws.Range("T46").Value = 0
$ws.Range("L48").Value = 0
$ws.Range("R48").Value = 4
$ws.Range("T48").Value = 4
$ws.Range("L49").Value = 0
$ws.Range("R49").Value = 3
$ws.Range("T49").Value = 2
$ws.Range("L50").Value = 0
$ws.Range("R50").Value = 46
$ws.Range("T50").Value = 0
$ws.Range("L52").Value = 0
$ws.Range("L53").Value = 0
$ws.Range("R55").Value = 1
$ws.Range("R56").Value = 1
$ws.Range("R57").Value = 1
$ws.Range("R58").Value = 2
$ws.Range("L62").Value = 0
$ws.Range("R62").Value = 2
$ws.Range("R63").Value = 1
$ws.Range("R64").Value = 2
$ws.Range("L66").Value = 0
$ws.Range("R66").Value = 4
$ws.Range("R67").Value = 2
$ws.Range("L69").Value = 0
$ws.Range("R69").Value = 4
$ws.Range("C72").Value = 120
$ws.Range("C83").Value = 0
